$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- F2: update the rpc-reply message-id UUID ---
$f2 = "<rpc-reply message-id=""urn:uuid:722f5cbf-9107-4beb-8084-6b3f9dedca69"">`n  <data/>`n</rpc-reply>`n"
$ws.Range("F2").Value = $f2

# --- G2: update protocol identifier/name entries ---
$g2 = @"
<edit-config>
    <target>
     <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
              <name>default</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>default</name>
              </config>
              <bgp>
                <global>
                  <config>
                    <as>65000</as>
                  </config>
                </global>
                <neighbors>
                  <neighbor>
                    <neighbor-address>192.168.1.2</neighbor-address>
                    <config>
                      <neighbor-address>192.168.1.2</neighbor-address>
                      <peer-as>65123</peer-as>
                      <description>BGP Neighbor 65123 Description</description>
                    </config>
                  </neighbor>
                </neighbors>
              </bgp>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
"@

$ws.Range("G2").Value = $g2
